$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in B2
$ws.Range("B2").Value = 175292

# Copy style from A2 (the existing label cell) to A3 and A4 for consistent formatting
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New row 3: granodiorite
$ws.Range("A3").Value = "granodiorite"
$ws.Range("B3").Value = 1713

# New row 4: quartz monzodiorite / quartz monzogabbro (embedded newline)
$ws.Range("A4").Value = "quartz monzodiorite`nquartz monzogabbro"
$ws.Range("B4").Value = 297
